$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- theta std-error column ("C") - refreshed bootstrap standard errors ---
$ws.Range("C3").Value  = "(0.13)"
$ws.Range("C5").Value  = "(0.21)"
$ws.Range("C7").Value  = "(0.15)"
$ws.Range("C9").Value  = "(0.01)"
$ws.Range("C11").Value = "(0.81)"
$ws.Range("C13").Value = "(0.26)"
$ws.Range("C15").Value = "(0.3)"

# --- lambda std-error column ("D") - refreshed bootstrap standard errors ---
$ws.Range("D3").Value  = "(0.02)"
$ws.Range("D5").Value  = "(0.05)"
$ws.Range("D7").Value  = "(0.04)"
$ws.Range("D9").Value  = "(0.01)"
$ws.Range("D11").Value = "(0.32)"
$ws.Range("D13").Value = "(0.01)"
$ws.Range("D15").Value = "(0.19)"

# --- proportion-drinking column ("E") - new values + std errors (bootstrapping) ---
$ws.Range("E2").Value  = 0.182455
$ws.Range("E3").Value  = "(0.00004)"
$ws.Range("E4").Value  = 0.167269
$ws.Range("E5").Value  = "(0.00002)"
$ws.Range("E6").Value  = 0.122975
$ws.Range("E7").Value  = "(0.00002)"
$ws.Range("E8").Value  = 0.101189
$ws.Range("E9").Value  = "(0.00000)"
$ws.Range("E10").Value = 0.092802
$ws.Range("E11").Value = "(0.00003)"
$ws.Range("E12").Value = 0.09722799999999999
$ws.Range("E13").Value = "(0.00000)"
$ws.Range("E14").Value = 0.089807
$ws.Range("E15").Value = "(0.00003)"
